$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($i = 11; $i -le 20; $i++) {
    $ws.Cells.Item($i, 1).Value = $i
}
